# Commit: "added final result calcs after adjusting ML labels"
# The vader-sentiment ML labels (column C, the buy/hold signal) were
# re-generated and a handful of dates flip between 0 and 1. Every other
# cell on the sheet (D/E running-balance columns and the H2/H3/I3/J3
# result-summary cells) is formula-driven off column C, so updating the
# raw label inputs and letting Excel recalculate reproduces the rest of
# the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Updated ML labels (column C) for the dates whose prediction changed
$changedLabels = @{
    3  = 0
    8  = 0
    9  = 1
    11 = 0
    20 = 0
    23 = 0
    27 = 1
    30 = 0
    31 = 1
    35 = 0
    41 = 0
    45 = 0
    46 = 0
    48 = 1
    49 = 1
    50 = 1
    51 = 1
}

foreach ($row in $changedLabels.Keys) {
    $ws.Cells.Item($row, 3).Value = $changedLabels[$row]
}

[void]$excel.Calculate()

# Matches the author's final cursor position recorded in the saved file
[void]$ws.Range("J3").Select()
